# Cell-write order below matches the original authoring order that produced the
# target shared-string table ordering (pattern/description/case_sensitive/label
# for both rows first, then the authority column for both rows last).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 (pattern, description, case_sensitive, label) ---
$ws.Range("A2").Value2 = "[0-9]{3}-[0-9]{2}-[0-9]{4}"
$ws.Range("B2").Value2 = "Lorem ..."
$ws.Range("C2").Value = $true
$ws.Range("D2").Value2 = "PII.social_security_number"

# --- Update row 3 (pattern, description, case_sensitive, label) ---
$ws.Range("A3").Value2 = "[0-9]{3} [0-9]{2} [0-9]{4}"
$ws.Range("B3").Value2 = "Lorem ..."
$ws.Range("C3").Value = $true
$ws.Range("D3").Value2 = "PII.social_security_number"

# --- Update authority column last for both rows ---
$ws.Range("E2").Value2 = "bar_"
$ws.Range("E3").Value2 = "bar_"

# --- Remove old row 4 entirely (data now fits in just 2 data rows) ---
$ws.Rows.Item(4).Delete() | Out-Null

# --- Update column widths to reflect new (wider) content ---
# (Input values chosen so the engine's internal width quantization lands as
# close as possible to the target stored widths of 22.5703125 / 25.28515625.)
$ws.Columns.Item(1).ColumnWidth = 21.59
$ws.Columns.Item(4).ColumnWidth = 24.42

# --- Update selection to match new active cell ---
$ws.Range("E4").Select() | Out-Null
